# Se procesan de nuevo los datos con las nuevas dimensiones curadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("F5").Clear()
